$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) "Stage 1: Time-series analysis" -> "Stage 1: Machine Learning"
# ------------------------------------------------------------------
$d.Content.Find.Execute("Stage 1: Time-series analysis", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Stage 1: Machine Learning", 2)

# ------------------------------------------------------------------
# 2) Rebuild the 3-column table into the new 6-column layout
# ------------------------------------------------------------------
$t = $d.Tables(1)

# --- grow 3 -> 6 columns -------------------------------------------------
$t.Columns.Add() | Out-Null
$t.Columns.Add() | Out-Null
$t.Columns.Add() | Out-Null

# --- column widths (dxa / 20 = points) -----------------------------------
$t.Columns(1).Width = 61.8    # 1236 dxa
$t.Columns(2).Width = 56.7    # 1134 dxa
$t.Columns(3).Width = 70.9    # 1418 dxa
$t.Columns(4).Width = 77.95   # 1559 dxa
$t.Columns(5).Width = 77.95   # 1559 dxa
$t.Columns(6).Width = 70.9    # 1418 dxa

# --- fixed layout + overall table width ----------------------------------
$t.AllowAutoFit = $false               # -> <w:tblLayout w:type="fixed"/>
$t.PreferredWidthType = 3              # wdPreferredWidthPoints
$t.PreferredWidth = 416.2              # 8324 dxa

# --- header row height ----------------------------------------------------
$t.Rows(1).Height = 39.2                # 784 dxa

# --- header row text --------------------------------------------------
# Column 1 ("Date of resignation") is unchanged.
$t.Rows(1).Cells(2).Range.Text = "Age"
$t.Rows(1).Cells(3).Range.Text = "Last promotion Date"
$t.Rows(1).Cells(4).Range.Text = "Increment percentage"
$t.Rows(1).Cells(5).Range.Text = "Dept."
$t.Rows(1).Cells(6).Range.Text = "Current salary"

# --- data rows --------------------------------------------------------
$t.Rows(2).Cells(1).Range.Text = "Anant"
$t.Rows(2).Cells(2).Range.Text = "54"
$t.Rows(2).Cells(3).Range.Text = "01/01/2023"
$t.Rows(2).Cells(4).Range.Text = "3%"
$t.Rows(2).Cells(5).Range.Text = "Test Analyst"
$t.Rows(2).Cells(6).Range.Text = "4.5 lpa"

$t.Rows(3).Cells(1).Range.Text = "Bazan"
$t.Rows(3).Cells(2).Range.Text = "36"
$t.Rows(3).Cells(3).Range.Text = "02/03/2025"
$t.Rows(3).Cells(4).Range.Text = "15%"
$t.Rows(3).Cells(5).Range.Text = "Programming"
$t.Rows(3).Cells(6).Range.Text = "8.5 lpa"

$t.Rows(4).Cells(1).Range.Text = "Chandra"
$t.Rows(4).Cells(2).Range.Text = "25"
$t.Rows(4).Cells(3).Range.Text = "02/03/205"
$t.Rows(4).Cells(4).Range.Text = "10%"
$t.Rows(4).Cells(5).Range.Text = "Analyst"
$t.Rows(4).Cells(6).Range.Text = "4.7 lpa"

# Row 5 stays blank (it already picked up the 3 new empty cells).
